$wb = $excel.ActiveWorkbook

# ----- ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H16").Value = 9504.5
$ws.Range("I16").Value = 9009
$ws.Range("J16").Value = 10000
$ws.Range("K16").Value = 9009
$ws.Range("L16").Value = 10000
$ws.Range("M16").Value = -8779
$ws.Range("N16").Value = -10460
$ws.Range("H62").Value = 333336670
$ws.Range("I62").Value = 333336670
$ws.Range("K62").Value = 333336670
$ws.Range("M62").Value = -333336046
$ws.Range("H65").Value = 333336670
$ws.Range("I65").Value = 333336670
$ws.Range("K65").Value = 1666683350
$ws.Range("M65").Value = -1666680230
$ws.Range("H76").Value = 4611.2
$ws.Range("I76").Value = 4631
$ws.Range("J76").Value = 4598
$ws.Range("K76").Value = 4631
$ws.Range("L76").Value = 4598
$ws.Range("M76").Value = -4316
$ws.Range("N76").Value = -5228
$ws.Range("H79").Value = 4611.2
$ws.Range("I79").Value = 4631
$ws.Range("J79").Value = 4598
$ws.Range("K79").Value = 4631
$ws.Range("L79").Value = 4598
$ws.Range("M79").Value = -3539
$ws.Range("N79").Value = -6782
$ws.Range("H133").Value = 89998.5
$ws.Range("J133").Value = 89998.5
$ws.Range("L133").Value = 89998.5
$ws.Range("N133").Value = -100118.5
$ws.Range("H137").Value = 12748.158
$ws.Range("I137").Value = 14395
$ws.Range("J137").Value = 3965
$ws.Range("K137").Value = 43185
$ws.Range("L137").Value = 11895
$ws.Range("M137").Value = -40635
$ws.Range("N137").Value = -16995

# ----- ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H16").Value = 100
$ws.Range("I16").Value = 100
$ws.Range("K16").Value = 100
$ws.Range("M16").Value = 187
$ws.Range("H32").Value = 6838.3276
$ws.Range("I32").Value = 6841.482
$ws.Range("K32").Value = 6841.482
$ws.Range("M32").Value = -6554.482
$ws.Range("H36").Value = 18800
$ws.Range("I36").Value = 18800
$ws.Range("K36").Value = 18800
$ws.Range("M36").Value = -18454
$ws.Range("H97").Value = 7696706.5
$ws.Range("I97").Value = 6369.5293
$ws.Range("K97").Value = 6369.5293
$ws.Range("M97").Value = -5873.5293

# ----- BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 12870.96
$ws.Range("I94").Value = 15693.556
$ws.Range("J94").Value = 5612.857
$ws.Range("K94").Value = 15693.556
$ws.Range("L94").Value = 5612.857
$ws.Range("M94").Value = -15242.556
$ws.Range("N94").Value = -6514.857

# ----- CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 15317.182
$ws.Range("I31").Value = 23252.5
$ws.Range("K31").Value = 23252.5
$ws.Range("M31").Value = -22957.5
$ws.Range("H34").Value = 15317.182
$ws.Range("I34").Value = 23252.5
$ws.Range("K34").Value = 23252.5
$ws.Range("M34").Value = -23050.5
$ws.Range("H114").Value = 62552.5
$ws.Range("J114").Value = 62552.5
$ws.Range("L114").Value = 62552.5
$ws.Range("N114").Value = -71230.5
$ws.Range("H141").Value = 326881.94
$ws.Range("J141").Value = 430635
$ws.Range("L141").Value = 430635
$ws.Range("N141").Value = -440995

# ----- CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H44").Value = 2922.5557
$ws.Range("I44").Value = 344.75
$ws.Range("J44").Value = 4984.8
$ws.Range("K44").Value = 1034.25
$ws.Range("L44").Value = 14954.4
$ws.Range("M44").Value = -636.25
$ws.Range("N44").Value = -15750.4
$ws.Range("H68").Value = 13140.637
$ws.Range("I68").Value = 2175
$ws.Range("J68").Value = 15577.444
$ws.Range("K68").Value = 6525
$ws.Range("L68").Value = 46732.33199999999
$ws.Range("M68").Value = -5714
$ws.Range("N68").Value = -48354.33199999999
$ws.Range("H69").Value = 3888.4443
$ws.Range("I69").Value = 2000
$ws.Range("J69").Value = 4124.5
$ws.Range("K69").Value = 6000
$ws.Range("L69").Value = 12373.5
$ws.Range("M69").Value = -5189
$ws.Range("N69").Value = -13995.5
$ws.Range("H71").Value = 13140.637
$ws.Range("I71").Value = 2175
$ws.Range("J71").Value = 15577.444
$ws.Range("K71").Value = 19575
$ws.Range("L71").Value = 140196.996
$ws.Range("M71").Value = -15519
$ws.Range("N71").Value = -148308.996
$ws.Range("H72").Value = 3888.4443
$ws.Range("I72").Value = 2000
$ws.Range("J72").Value = 4124.5
$ws.Range("K72").Value = 18000
$ws.Range("L72").Value = 37120.5
$ws.Range("M72").Value = -13944
$ws.Range("N72").Value = -45232.5
$ws.Range("H80").Value = 87199.8
$ws.Range("I80").Value = 7595
$ws.Range("J80").Value = 107101
$ws.Range("K80").Value = 22785
$ws.Range("L80").Value = 321303
$ws.Range("M80").Value = -21849
$ws.Range("N80").Value = -323175
$ws.Range("H83").Value = 87199.8
$ws.Range("I83").Value = 7595
$ws.Range("J83").Value = 107101
$ws.Range("K83").Value = 68355
$ws.Range("L83").Value = 963909
$ws.Range("M83").Value = -63675
$ws.Range("N83").Value = -973269
$ws.Range("H86").Value = 406.75
$ws.Range("J86").Value = 393
$ws.Range("L86").Value = 1179
$ws.Range("N86").Value = -3551
$ws.Range("H87").Value = 16222.846
$ws.Range("I87").Value = 11506.75
$ws.Range("J87").Value = 18318.889
$ws.Range("K87").Value = 34520.25
$ws.Range("L87").Value = 54956.667
$ws.Range("M87").Value = -33272.25
$ws.Range("N87").Value = -57452.667
$ws.Range("H89").Value = 406.75
$ws.Range("J89").Value = 393
$ws.Range("L89").Value = 3537
$ws.Range("N89").Value = -15393
$ws.Range("H90").Value = 16222.846
$ws.Range("I90").Value = 11506.75
$ws.Range("J90").Value = 18318.889
$ws.Range("K90").Value = 103560.75
$ws.Range("L90").Value = 164870.001
$ws.Range("M90").Value = -97320.75
$ws.Range("N90").Value = -177350.001
$ws.Range("H112").Value = 4749.8335
$ws.Range("I112").Value = 3000
$ws.Range("J112").Value = 6499.6665
$ws.Range("K112").Value = 9000
$ws.Range("L112").Value = 19498.9995
$ws.Range("M112").Value = -7892
$ws.Range("N112").Value = -21714.9995

# ----- GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 8562.947
$ws.Range("J80").Value = 6179.5713
$ws.Range("L80").Value = 6179.5713
$ws.Range("N80").Value = -8175.5713
$ws.Range("H83").Value = 8562.947
$ws.Range("J83").Value = 6179.5713
$ws.Range("L83").Value = 30897.8565
$ws.Range("N83").Value = -40881.85649999999
$ws.Range("H97").Value = 7508
$ws.Range("I97").Value = 8560.0625
$ws.Range("K97").Value = 8560.0625
$ws.Range("M97").Value = -8064.0625

# ----- LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 20802.031
$ws.Range("I40").Value = 26075.555
$ws.Range("K40").Value = 26075.555
$ws.Range("M40").Value = -25939.555
$ws.Range("H122").Value = 5136.1904
$ws.Range("I122").Value = 5294.5
$ws.Range("J122").Value = 4186.3335
$ws.Range("K122").Value = 15883.5
$ws.Range("L122").Value = 12559.0005
$ws.Range("M122").Value = -13433.5
$ws.Range("N122").Value = -17459.0005

# ----- WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H69").Value = 30450
$ws.Range("J69").Value = 30450
$ws.Range("L69").Value = 30450
$ws.Range("N69").Value = -31948
$ws.Range("H72").Value = 30450
$ws.Range("J72").Value = 30450
$ws.Range("L72").Value = 91350
$ws.Range("N72").Value = -98838
$ws.Range("H132").Value = 10003.9
$ws.Range("I132").Value = 13230.971
$ws.Range("K132").Value = 39692.913
$ws.Range("M132").Value = -37162.913
